# [BAU_FIX] impact enum trailing space fix
#
# The "Impact" enum values used in the RiskRegister sheet ("4 - Significant
# impact" and "5 - Major impact") were missing a trailing space compared to
# the rest of the enum list. Fix the typo by adding the trailing space back,
# which also updates every cell in the RiskRegister sheet that currently
# holds one of those two values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RiskRegister")

$rng = $ws.UsedRange

# Order matters for shared-string table layout: fix "4 - Significant impact"
# before "5 - Major impact".
$rng.Replace("4 - Significant impact", "4 - Significant impact ")
$rng.Replace("5 - Major impact", "5 - Major impact ")

# Restore the active selection on the RiskRegister sheet.
$ws.Activate()
$ws.Range("L11").Select()
